$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new time-tracking entry row (row 38). Only column A carries a
# non-default (date) style, so copy just that cell's formatting from the
# row above it — columns B and C already use their column's default style.
$ws.Range("A37").Copy() | Out-Null
$ws.Range("A38").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A38").Value = 45636
$ws.Range("B38").Value = 1.5
$ws.Range("C38").Value = "Created Node server, ported Flask server functionality to Node"

# Keep the selection in sync with where Excel would leave the cursor after data entry
$ws.Range("A39").Select() | Out-Null

$wb.Application.Calculate() | Out-Null
